$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200886964797974
$ws.Range("B1").Value = 3.780808448791504
$ws.Range("C1").Value = 3.264467239379883
$ws.Range("D1").Value = 2.598045825958252
$ws.Range("E1").Value = 1.2761470079422
